$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to the batter's name.
$ws.Name = "Fabian Allen"

# Insert a new column A ("matchNo") before the existing data, shifting
# everything (teamName..result) one column to the right (B..M).
$ws.Columns("A:A").Insert()

# ---- Row 1 (headers) ----
$ws.Cells.Item(1, 1).Value = "matchNo"

# ---- Row 2 (existing match row, now with matchNo filled in) ----
$ws.Cells.Item(2, 1).Value = "45th"

# ---- Row 3 (new match row) ----
$ws.Cells.Item(3, 1).Value = "32nd"
$ws.Cells.Item(3, 2).Value = "Punjab Kings"
$ws.Cells.Item(3, 3).Value = "Fabian Allen"
$ws.Cells.Item(3, 4).Value = "'"
$ws.Cells.Item(3, 5).Value = "'0"
$ws.Cells.Item(3, 6).Value = "'1"
$ws.Cells.Item(3, 7).Value = "'0"
$ws.Cells.Item(3, 8).Value = "'0"
$ws.Cells.Item(3, 9).Value = "'0.00"
$ws.Cells.Item(3, 10).Value = "Rajasthan Royals"
$ws.Cells.Item(3, 11).Value = "Dubai (DSC)"
$ws.Cells.Item(3, 12).Value = "September 21"
$ws.Cells.Item(3, 13).Value = "Royals won by 2 runs"

# ---- Row 4 (new match row) ----
$ws.Cells.Item(4, 1).Value = "14th"
$ws.Cells.Item(4, 2).Value = "Punjab Kings"
$ws.Cells.Item(4, 3).Value = "Fabian Allen"
$ws.Cells.Item(4, 4).Value = "c Warner b Ahmed"
$ws.Cells.Item(4, 5).Value = "'6"
$ws.Cells.Item(4, 6).Value = "'11"
$ws.Cells.Item(4, 7).Value = "'0"
$ws.Cells.Item(4, 8).Value = "'0"
$ws.Cells.Item(4, 9).Value = "'54.54"
$ws.Cells.Item(4, 10).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(4, 11).Value = "Chennai"
$ws.Cells.Item(4, 12).Value = "April 21"
$ws.Cells.Item(4, 13).Value = "Sunrisers won by 9 wickets (with 8 balls remaining)"

# The leading single-quotes above force Excel to keep numeric-looking
# literals ("0", "1", "0.00", "6", "11", "54.54") and the blank "states"
# cells as TEXT (matching every other cell in this sheet), exactly like
# the rest of the scraped data. Drop the resulting quote-prefix styling
# so the cells end up with the sheet's plain default format.
$ws.Range("D3:I4").ClearFormats()
